# Applies the text edits to Banking_EDA_Skeleton.pptx:
#  - Slide 2 ("Hypotheses"): reword the Null/Alternative hypothesis bullets
#  - Slide 3 ("The data"):   reword the "All of them" bullet
#  - Slide 5 ("Wrap-up"):    merge the split "Austin Youngblood" contact runs
#
# Note: for each edit we first stomp the target paragraph with a throwaway
# value before writing the real text. The host's TextRange.Text setter
# otherwise tries to preserve as much of the previous string as possible
# (a longest-common-substring style diff) and ends up fragmenting the
# paragraph into several <a:r> runs instead of the single run PowerPoint
# itself would produce when a user retypes/selects-and-replaces the text.

$p = $ppt.ActivePresentation

# --- Slide 2 ("Hypotheses") : Content Placeholder 2 ---------------------
$slide2 = $p.Slides.Item(2)
$shape2 = $slide2.Shapes.Item(2)
$tr2 = $shape2.TextFrame.TextRange

$nullHyp = $tr2.Paragraphs(1, 1)
$nullHyp.Text = "X"
$nullHyp.Text = "Null Hypothesis (H0): There is no significant difference the effect of loan length has on the interest rate."

$altHyp = $tr2.Paragraphs(4, 1)
$altHyp.Text = "X"
$altHyp.Text = "Alternative Hypothesis (H1): There is a significant difference the effect of loan length has on interest rates."

# --- Slide 3 ("The data") : Content Placeholder 2 ------------------------
$slide3 = $p.Slides.Item(3)
$shape3 = $slide3.Shapes.Item(2)
$tr3 = $shape3.TextFrame.TextRange

$featureBullet = $tr3.Paragraphs(2, 1)
$featureBullet.Text = "X"
$featureBullet.Text = "Interest rate as the target and Loan Length as the feature"

# --- Slide 5 ("Wrap-up") : Content Placeholder 2 --------------------------
$slide5 = $p.Slides.Item(5)
$shape5 = $slide5.Shapes.Item(2)
$tr5 = $shape5.TextFrame.TextRange

$austinLine = $tr5.Paragraphs(3, 1)
$austinLine.Text = "X"
$austinLine.Text = "Austin Youngblood - Austin.m.youngblood@gmail.com"
